$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.560.92'
$ws.Range('E2').Value = '  +0.21%  '
$ws.Range('D3').Value = '1.913.45'
$ws.Range('E3').Value = '  -0.01%  '
$c = $ws.Range('D4')
$c.NumberFormat = "@"
$c.Value = '1.006'
$c.Style = "Normal"
$ws.Range('E4').Value = '  +0.57%  '
$c = $ws.Range('D5')
$c.NumberFormat = "@"
$c.Value = '325.84'
$c.Style = "Normal"
$ws.Range('E5').Value = '  -0.28%  '
$c = $ws.Range('D6')
$c.NumberFormat = "@"
$c.Value = '1.005'
$c.Style = "Normal"
$ws.Range('E6').Value = '  +0.46%  '
$ws.Range('E7').Value = '  +1.32%  '
$c = $ws.Range('D8')
$c.NumberFormat = "@"
$c.Value = '0.4078'
$c.Style = "Normal"
$ws.Range('E8').Value = '  -0.46%  '
$c = $ws.Range('D9')
$c.NumberFormat = "@"
$c.Value = '0.08167'
$c.Style = "Normal"
$ws.Range('E9').Value = '  +1.67%  '
$c = $ws.Range('D10')
$c.NumberFormat = "@"
$c.Value = '1.013'
$c.Style = "Normal"
$ws.Range('E10').Value = '  +0.28%  '
$c = $ws.Range('D11')
$c.NumberFormat = "@"
$c.Value = '23.49'
$c.Style = "Normal"
$ws.Range('E11').Value = '  +5.00%  '
$ws.Range('D12').Value = '1.921.78'
$ws.Range('E12').Value = '  +0.53%  '
$c = $ws.Range('D13')
$c.NumberFormat = "@"
$c.Value = '6.026'
$c.Style = "Normal"
$ws.Range('E13').Value = '  +1.35%  '
$c = $ws.Range('D14')
$c.NumberFormat = "@"
$c.Value = '7.121'
$c.Style = "Normal"
$ws.Range('E14').Value = '  -0.47%  '
$c = $ws.Range('D15')
$c.NumberFormat = "@"
$c.Value = '90.56'
$c.Style = "Normal"
$ws.Range('E15').Value = '  +1.17%  '
$c = $ws.Range('D16')
$c.NumberFormat = "@"
$c.Value = '0.06794'
$c.Style = "Normal"
$ws.Range('E16').Value = '  +2.87%  '
$c = $ws.Range('D17')
$c.NumberFormat = "@"
$c.Value = '1.007'
$c.Style = "Normal"
$ws.Range('E17').Value = '  +0.56%  '
$c = $ws.Range('D18')
$c.NumberFormat = "@"
$c.Value = '0.00001043'
$c.Style = "Normal"
$ws.Range('E18').Value = '  +1.10%  '
$c = $ws.Range('D19')
$c.NumberFormat = "@"
$c.Value = '17.74'
$c.Style = "Normal"
$ws.Range('E19').Value = '  -0.16%  '
$c = $ws.Range('D20')
$c.NumberFormat = "@"
$c.Value = '1.005'
$c.Style = "Normal"
$ws.Range('E20').Value = '  +0.47%  '
$ws.Range('D21').Value = '29.565.58'
$ws.Range('E21').Value = '  +0.23%  '
$c = $ws.Range('D22')
$c.NumberFormat = "@"
$c.Value = '5.617'
$c.Style = "Normal"
$ws.Range('E22').Value = '  +1.22%  '
$c = $ws.Range('D23')
$c.NumberFormat = "@"
$c.Value = '11.80'
$c.Style = "Normal"
$ws.Range('E23').Value = '  +2.48%  '
$c = $ws.Range('D24')
$c.NumberFormat = "@"
$c.Value = '2.171'
$c.Style = "Normal"
$ws.Range('E24').Value = '  -1.49%  '
$ws.Range('D25').Value = '2.161.99'
$ws.Range('E25').Value = '  +0.96%  '
$c = $ws.Range('D26')
$c.NumberFormat = "@"
$c.Value = '154.58'
$c.Style = "Normal"
$ws.Range('E26').Value = '  +0.64%  '
$c = $ws.Range('D27')
$c.NumberFormat = "@"
$c.Value = '20.11'
$c.Style = "Normal"
$ws.Range('E27').Value = '  +1.54%  '
$c = $ws.Range('D28')
$c.NumberFormat = "@"
$c.Value = '6.339'
$c.Style = "Normal"
$ws.Range('E28').Value = '  +8.85%  '
$c = $ws.Range('D29')
$c.NumberFormat = "@"
$c.Value = '2.108'
$c.Style = "Normal"
$ws.Range('E29').Value = '  -1.15%  '
$c = $ws.Range('D30')
$c.NumberFormat = "@"
$c.Value = '119.66'
$c.Style = "Normal"
$ws.Range('E30').Value = '  +1.76%  '
$c = $ws.Range('D31')
$c.NumberFormat = "@"
$c.Value = '1.030'
$c.Style = "Normal"
$ws.Range('E31').Value = '  -2.83%  '
$c = $ws.Range('D32')
$c.NumberFormat = "@"
$c.Value = '0.09582'
$c.Style = "Normal"
$ws.Range('E32').Value = '  +0.32%  '
$c = $ws.Range('D33')
$c.NumberFormat = "@"
$c.Value = '5.555'
$c.Style = "Normal"
$ws.Range('E33').Value = '  +3.15%  '
$c = $ws.Range('D34')
$c.NumberFormat = "@"
$c.Value = '1.397'
$c.Style = "Normal"
$ws.Range('E34').Value = '  -1.67%  '
$c = $ws.Range('D35')
$c.NumberFormat = "@"
$c.Value = '3.555'
$c.Style = "Normal"
$ws.Range('E35').Value = '  -0.52%  '
$c = $ws.Range('D36')
$c.NumberFormat = "@"
$c.Value = '0.02274'
$c.Style = "Normal"
$ws.Range('E36').Value = '  +0.88%  '
$c = $ws.Range('D37')
$c.NumberFormat = "@"
$c.Value = '0.06118'
$c.Style = "Normal"
$ws.Range('E37').Value = '  +0.31%  '
$c = $ws.Range('D38')
$c.NumberFormat = "@"
$c.Value = '1.181'
$c.Style = "Normal"
$ws.Range('E38').Value = '  +0.58%  '
$ws.Range('B39').Value = 'TheSandbox'
$ws.Range('C39').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$c = $ws.Range('D39')
$c.NumberFormat = "@"
$c.Value = '0.5948'
$c.Style = "Normal"
$ws.Range('E39').Value = '  +1.22%  '
$ws.Range('B40').Value = 'Aptos'
$ws.Range('C40').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$c = $ws.Range('D40')
$c.NumberFormat = "@"
$c.Value = '10.79'
$c.Style = "Normal"
$ws.Range('E40').Value = '  +6.66%  '
$c = $ws.Range('D41')
$c.NumberFormat = "@"
$c.Value = '7.959'
$c.Style = "Normal"
$ws.Range('E41').Value = '  -4.70%  '
$c = $ws.Range('D42')
$c.NumberFormat = "@"
$c.Value = '0.1855'
$c.Style = "Normal"
$ws.Range('E42').Value = '  +0.69%  '
$c = $ws.Range('D43')
$c.NumberFormat = "@"
$c.Value = '2.461'
$c.Style = "Normal"
$ws.Range('E43').Value = '  +0.69%  '
$c = $ws.Range('D44')
$c.NumberFormat = "@"
$c.Value = '1.281'
$c.Style = "Normal"
$ws.Range('E44').Value = '  -0.98%  '
$c = $ws.Range('D45')
$c.NumberFormat = "@"
$c.Value = '0.07730'
$c.Style = "Normal"
$ws.Range('E45').Value = '  -3.67%  '
$c = $ws.Range('D46')
$c.NumberFormat = "@"
$c.Value = '12.44'
$c.Style = "Normal"
$ws.Range('E46').Value = '  +2.29%  '
$c = $ws.Range('D47')
$c.NumberFormat = "@"
$c.Value = '0.5582'
$c.Style = "Normal"
$ws.Range('E47').Value = '  +0.75%  '
$c = $ws.Range('D48')
$c.NumberFormat = "@"
$c.Value = '1.955'
$c.Style = "Normal"
$ws.Range('E48').Value = '  +1.41%  '
$c = $ws.Range('D49')
$c.NumberFormat = "@"
$c.Value = '115.14'
$c.Style = "Normal"
$ws.Range('E49').Value = '  +1.32%  '
$c = $ws.Range('D50')
$c.NumberFormat = "@"
$c.Value = '72.89'
$c.Style = "Normal"
$ws.Range('E50').Value = '  +1.81%  '
$c = $ws.Range('D51')
$c.NumberFormat = "@"
$c.Value = '1.054'
$c.Style = "Normal"
$ws.Range('E51').Value = '  +2.14%  '
